$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reading-data")
$ws.Activate() | Out-Null

# "pages" (F) / "book" (G) reading-log entries for 2019-02-03 .. 2019-03-04
# (rows 41-67). Row 41 already existed with a placeholder pages count that
# gets corrected; rows 42-67 are newly logged days.
$pages = @{
    41 = 27
    42 = 19
    43 = 15
    44 = 17
    45 = 0
    46 = 22
    47 = 0
    48 = 7
    49 = 13
    50 = 19
    51 = 11
    52 = 14
    53 = 12
    54 = 8
    55 = 0
    56 = 18
    57 = 17
    58 = 0
    59 = 4
    60 = 9
    61 = 20
    62 = 5
    63 = 36
    64 = 0
    65 = 0
    66 = 8
    67 = 20
}

foreach ($row in 41..67) {
    $ws.Cells.Item($row, 6).Value = $pages[$row]
    $ws.Cells.Item($row, 7).Value = 3
}

# Leave the view scrolled/selected where the author left off editing.
$ws.Range("I64").Select() | Out-Null
